# CO2 conférences.xlsx - "Collectif" sheet refinements
# - bump duration / participants inputs
# - adjust several % breakdown inputs (Voyages + Repas sections)
# - add "Nombre max" row + two new K-column (Nombre/trajet) helper formulas
# - move active selection to K9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collectif")

# --- Informations générales ---------------------------------------------
$ws.Range("D4").Value = 6    # Durée (jours) : 5 -> 6
$ws.Range("D5").Value = 76   # Participants : 100 -> 76

# --- Voyages : répartition (% des participants) -------------------------
$ws.Range("C10").Value = 0.15
$ws.Range("C11").Value = 0.07
$ws.Range("C12").Value = 0.05
$ws.Range("C13").Value = 0.06
$ws.Range("C14").Value = 0.67

# New "Nombre / trajet" helper column for the first two rows of the table
$ws.Range("K8").Formula = "=E8*F8*G8"
$ws.Range("K9").Formula = "=E9*F9*G9"

# --- Repas : répartition (% des participants) ----------------------------
$ws.Range("C20").Value = 0.48
$ws.Range("C21").Value = 0.08
$ws.Range("C22").Value = 0.45

# --- New label row between the "Voyages" total and the "Repas" header ----
$ws.Range("E19").Value = "Nombre max"

# --- Restore the active selection shown in the saved workbook -----------
$ws.Range("K9").Select()
